$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{Row=2; I='aa'; J='Agree/Accept'}
    @{Row=13; I='aa'; J='Agree/Accept'}
    @{Row=15; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=28; I='aa'; J='Agree/Accept'}
    @{Row=42; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=55; I='sv'; J='Statement-opinion'}
    @{Row=61; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=64; I='sd'; J='Statement-non-opinion'}
    @{Row=73; I='sd'; J='Statement-non-opinion'}
    @{Row=83; I='sv'; J='Statement-opinion'}
    @{Row=84; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=90; I='sv'; J='Statement-opinion'}
    @{Row=92; I='sv'; J='Statement-opinion'}
    @{Row=94; I='sv'; J='Statement-opinion'}
    @{Row=110; I='sv'; J='Statement-opinion'}
    @{Row=141; I='sv'; J='Statement-opinion'}
    @{Row=142; I='sv'; J='Statement-opinion'}
    @{Row=147; I='sv'; J='Statement-opinion'}
    @{Row=149; I='sd'; J='Statement-non-opinion'}
    @{Row=153; I='aa'; J='Agree/Accept'}
    @{Row=160; I='sv'; J='Statement-opinion'}
    @{Row=165; I='sv'; J='Statement-opinion'}
    @{Row=196; I='sd'; J='Statement-non-opinion'}
    @{Row=198; I='sv'; J='Statement-opinion'}
    @{Row=204; I='sd'; J='Statement-non-opinion'}
    @{Row=214; I='sv'; J='Statement-opinion'}
    @{Row=217; I='aa'; J='Agree/Accept'}
    @{Row=269; I='aa'; J='Agree/Accept'}
    @{Row=273; I='sd'; J='Statement-non-opinion'}
    @{Row=274; I='aa'; J='Agree/Accept'}
    @{Row=289; I='ba'; J='Appreciation'}
    @{Row=334; I='sd'; J='Statement-non-opinion'}
    @{Row=340; I='ba'; J='Appreciation'}
    @{Row=378; I='%'; J='Uninterpretable'}
    @{Row=385; I='sv'; J='Statement-opinion'}
    @{Row=387; I='sd'; J='Statement-non-opinion'}
    @{Row=390; I='sd'; J='Statement-non-opinion'}
    @{Row=405; I='sd'; J='Statement-non-opinion'}
    @{Row=409; I='sv'; J='Statement-opinion'}
    @{Row=410; I='sd'; J='Statement-non-opinion'}
    @{Row=419; I='sd'; J='Statement-non-opinion'}
    @{Row=421; I='sv'; J='Statement-opinion'}
    @{Row=429; I='aa'; J='Agree/Accept'}
    @{Row=446; I='ba'; J='Appreciation'}
    @{Row=451; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=456; I='sv'; J='Statement-opinion'}
    @{Row=462; I='sv'; J='Statement-opinion'}
    @{Row=467; I='sd'; J='Statement-non-opinion'}
    @{Row=471; I='sv'; J='Statement-opinion'}
    @{Row=472; I='sd'; J='Statement-non-opinion'}
    @{Row=484; I='sd'; J='Statement-non-opinion'}
    @{Row=491; I='sd'; J='Statement-non-opinion'}
    @{Row=499; I='sv'; J='Statement-opinion'}
    @{Row=509; I='sv'; J='Statement-opinion'}
    @{Row=511; I='ba'; J='Appreciation'}
    @{Row=520; I='sv'; J='Statement-opinion'}
    @{Row=521; I='ba'; J='Appreciation'}
    @{Row=537; I='%'; J='Uninterpretable'}
    @{Row=538; I='%'; J='Uninterpretable'}
    @{Row=544; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=552; I='aa'; J='Agree/Accept'}
    @{Row=562; I='sd'; J='Statement-non-opinion'}
    @{Row=575; I='ba'; J='Appreciation'}
    @{Row=580; I='sv'; J='Statement-opinion'}
    @{Row=581; I='ba'; J='Appreciation'}
    @{Row=582; I='sd'; J='Statement-non-opinion'}
    @{Row=587; I='sd'; J='Statement-non-opinion'}
    @{Row=596; I='ba'; J='Appreciation'}
    @{Row=600; I='ba'; J='Appreciation'}
    @{Row=611; I='sv'; J='Statement-opinion'}
    @{Row=612; I='ba'; J='Appreciation'}
    @{Row=614; I='sv'; J='Statement-opinion'}
    @{Row=619; I='ba'; J='Appreciation'}
    @{Row=631; I='aa'; J='Agree/Accept'}
    @{Row=632; I='ba'; J='Appreciation'}
    @{Row=651; I='sd'; J='Statement-non-opinion'}
    @{Row=674; I='sd'; J='Statement-non-opinion'}
    @{Row=681; I='aa'; J='Agree/Accept'}
    @{Row=698; I='aa'; J='Agree/Accept'}
    @{Row=704; I='sd'; J='Statement-non-opinion'}
    @{Row=730; I='sv'; J='Statement-opinion'}
    @{Row=739; I='sd'; J='Statement-non-opinion'}
    @{Row=741; I='sv'; J='Statement-opinion'}
    @{Row=748; I='aa'; J='Agree/Accept'}
    @{Row=749; I='sd'; J='Statement-non-opinion'}
    @{Row=752; I='sd'; J='Statement-non-opinion'}
    @{Row=755; I='sd'; J='Statement-non-opinion'}
    @{Row=758; I='sd'; J='Statement-non-opinion'}
    @{Row=761; I='sd'; J='Statement-non-opinion'}
    @{Row=764; I='sd'; J='Statement-non-opinion'}
    @{Row=781; I='aa'; J='Agree/Accept'}
    @{Row=782; I='sd'; J='Statement-non-opinion'}
    @{Row=785; I='aa'; J='Agree/Accept'}
    @{Row=797; I='aa'; J='Agree/Accept'}
    @{Row=799; I='sv'; J='Statement-opinion'}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
